# CV Discount Check Master File — update Accessories price for the
# "BOL CAMPER" model family (Sheet1 rows 34-38) from 15000 to 12000.
# Dependent formulas (K = SUM(C:J), L = K-H) on Sheet1, and the
# Report sheet's VLOOKUP-driven D12:D14 cells, recalculate automatically.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Sheet1")

$data.Range("J34:J38").Value = 12000
